$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2-28 hold a date ("Förändrad") stored as serial 45445 (2024-06-02).
# Update it to serial 45446 (2024-06-03) for every row, keeping existing date formatting.
for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45445) {
        $cell.Value = 45446
    }
}
